$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-6, columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
# This reflects the weekly re-shuffle of rows described in the commit message.

$rows = @(
    @{ Row = 2; D = 44176; L = "Primera"; M = 300; N = 5000; O = 6000; P = 5500; S = 3667 },
    @{ Row = 3; D = 44162; L = "Primera"; M = 100; N = 7000; O = 7000; P = 7000; S = 4667 },
    @{ Row = 4; D = 44162; L = "Segunda"; M = 100; N = 6500; O = 6500; P = 6500; S = 4333 },
    @{ Row = 5; D = 44159; L = "Segunda"; M = 200; N = 6500; O = 7000; P = 6750; S = 4500 },
    @{ Row = 6; D = 44169; L = "Primera"; M = 400; N = 5500; O = 6000; P = 5750; S = 3833 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L - Calidad
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M - Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N - Precio minimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O - Precio maximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r.Row, 19).Value = $r.S   # S - Precio $/Kg
}
